$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timesheet data entry: Wednesday (row 15) clocked in 15:00, out 18:00 ---
$ws.Range("C15").Value = 0.625
$ws.Range("D15").Value = 0.75

# --- Move the selection cursor down to D16 (next row's "Out" cell) ---
$ws.Range("D16").Select()

# --- Add another generated Print_Area defined name (mirrors the repeated
#     re-save artifact already present in the workbook) ---
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# --- Minor column-width nudge (columns were trimmed slightly) ---
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 8.8
$ws.Columns.Item(8).ColumnWidth = 9.8
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 8.8
$ws.Range("L1:AMK1").EntireColumn.ColumnWidth = 6
